$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B5").Value = 0.00222625182058704
$ws.Range("C5").Value = 0.001569787117169649
$ws.Range("D5").Value = 0.001325577226074754
$ws.Range("E5").Value = 0.0018572632978132
$ws.Range("F5").Value = 0.001738689689759812
$ws.Range("G5").Value = 0.00157962338536876

$ws.Range("B6").Value = 0.001991614944459785
$ws.Range("C6").Value = 0.001291350478748593
$ws.Range("D6").Value = 0.001244958982792127
$ws.Range("E6").Value = 0.001218624070470127
$ws.Range("F6").Value = 0.001429028797842887
$ws.Range("G6").Value = 0.001122456285256862

$ws.Range("B7").Value = 0.004350634231105001
$ws.Range("C7").Value = 0.004538637975386308
$ws.Range("D7").Value = 0.003314520384202478
$ws.Range("E7").Value = 0.004140273757784414
$ws.Range("F7").Value = 0.00426150058787862
$ws.Range("G7").Value = 0.003246756277608157

$ws.Range("B8").Value = 0.0005953469904289718
$ws.Range("C8").Value = 0.0004268645716496332
$ws.Range("D8").Value = 0.0005533050419349652
$ws.Range("E8").Value = 0.0007503747311340246
$ws.Range("F8").Value = 0.0007675501956847752
$ws.Range("G8").Value = 0.0006647402874995507
